$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

# Update the unit/description for FOM parameter (row 9)
$ws.Range("B9").Value = "EUR/(MW/(hours per model year))/8760"

# Add new row 15 for E2H parameter
$ws.Range("A15").Value = "E2H"
$ws.Range("B15").Value = "Coefficient (negative for heat pumps, positive for backpressure)"
